$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $NewValue
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "57.467.76"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "3.101.49"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "523.10"
$ws.Range("E5").Value = "  +1.03%  "

Set-TextValue "D6" "141.07"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.099.68"
$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "3.634.46"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("E14").Value = "  +1.11%  "

Set-TextValue "D15" "26.12"
$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").Value = "57.550.43"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "3.102.39"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  -0.71%  "

Set-TextValue "D22" "336.13"

$ws.Range("E23").Value = "  -0.04%  "

Set-TextValue "D24" "0.511"
$ws.Range("E24").Value = "  +2.44%  "

$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("E26").Value = "  -0.70%  "

$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").Value = "0.0₃0919"
$ws.Range("E28").Value = "  +1.70%  "

Set-TextValue "D29" "6.51"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +2.16%  "

Set-TextValue "D35" "156.90"
$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("E36").Value = "  +3.06%  "

Set-TextValue "D37" "6.11"
$ws.Range("E37").Value = "  +2.84%  "

Set-TextValue "D38" "27.06"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("E39").Value = "  +1.35%  "

Set-TextValue "D40" "0.0660"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").Value = "3.141.14"
$ws.Range("E41").Value = "  +1.12%  "

Set-TextValue "D42" "3.94"
$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("E43").Value = "  +4.49%  "

$ws.Range("E44").Value = "  +11.04%  "

$ws.Range("E45").Value = "  +0.66%  "

Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "2.296.94"
$ws.Range("E47").Value = "  +1.75%  "

Set-TextValue "D48" "0.0259"
$ws.Range("E48").Value = "  +0.26%  "

Set-TextValue "D49" "0.976"
$ws.Range("E49").Value = "  +4.92%  "

Set-TextValue "D50" "20.75"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("E51").Value = "  +2.21%  "

# Row 33/34 swap: EthereumClassic <-> Fetch.AI
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D33" "1.20"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "20.89"
$ws.Range("E34").Value = "  +0.85%  "
